{"js": "// Replace each two-digit multiplication problem's text with its new value.\n// Every `<w:t>` run in the practice-sheet table holds a unique \"AA\u00d7BB=\"\n// string, so a scoped, case-sensitive `body.search()` + `insertText(...,\n// \"Replace\")` per pair is safe (no cross-matches between old/new values).\nconst replacements = [\n  [\"46\u00d738=\", \"94\u00d768=\"],\n  [\"76\u00d740=\", \"88\u00d798=\"],\n  [\"89\u00d742=\", \"32\u00d742=\"],\n  [\"55\u00d755=\", \"97\u00d746=\"],\n  [\"14\u00d753=\", \"78\u00d744=\"],\n  [\"88\u00d717=\", \"27\u00d732=\"],\n  [\"27\u00d743=\", \"71\u00d798=\"],\n  [\"36\u00d770=\", \"62\u00d771=\"],\n  [\"78\u00d730=\", \"74\u00d767=\"],\n  [\"37\u00d720=\", \"37\u00d776=\"],\n  [\"49\u00d741=\", \"48\u00d796=\"],\n  [\"70\u00d743=\", \"79\u00d776=\"],\n  [\"82\u00d721=\", \"67\u00d772=\"],\n  [\"75\u00d731=\", \"55\u00d725=\"],\n  [\"94\u00d774=\", \"77\u00d746=\"],\n  [\"97\u00d761=\", \"22\u00d796=\"],\n  [\"62\u00d786=\", \"60\u00d766=\"],\n  [\"68\u00d779=\", \"74\u00d763=\"],\n  [\"50\u00d775=\", \"19\u00d740=\"],\n  [\"37\u00d736=\", \"61\u00d745=\"],\n  [\"83\u00d721=\", \"18\u00d750=\"],\n  [\"16\u00d743=\", \"40\u00d728=\"],\n  [\"56\u00d748=\", \"25\u00d746=\"],\n  [\"16\u00d748=\", \"82\u00d720=\"],\n  [\"39\u00d730=\", \"41\u00d797=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update every two-digit multiplication problem in the practice-sheet\n# table to its new value. Each \"AA\u00d7BB=\" string is unique in the document,\n# so a plain Find/Replace (one pair at a time, match case, whole story)\n# is safe and idempotent.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n    @(\"46\u00d738=\", \"94\u00d768=\"),\n    @(\"76\u00d740=\", \"88\u00d798=\"),\n    @(\"89\u00d742=\", \"32\u00d742=\"),\n    @(\"55\u00d755=\", \"97\u00d746=\"),\n    @(\"14\u00d753=\", \"78\u00d744=\"),\n    @(\"88\u00d717=\", \"27\u00d732=\"),\n    @(\"27\u00d743=\", \"71\u00d798=\"),\n    @(\"36\u00d770=\", \"62\u00d771=\"),\n    @(\"78\u00d730=\", \"74\u00d767=\"),\n    @(\"37\u00d720=\", \"37\u00d776=\"),\n    @(\"49\u00d741=\", \"48\u00d796=\"),\n    @(\"70\u00d743=\", \"79\u00d776=\"),\n    @(\"82\u00d721=\", \"67\u00d772=\"),\n    @(\"75\u00d731=\", \"55\u00d725=\"),\n    @(\"94\u00d774=\", \"77\u00d746=\"),\n    @(\"97\u00d761=\", \"22\u00d796=\"),\n    @(\"62\u00d786=\", \"60\u00d766=\"),\n    @(\"68\u00d779=\", \"74\u00d763=\"),\n    @(\"50\u00d775=\", \"19\u00d740=\"),\n    @(\"37\u00d736=\", \"61\u00d745=\"),\n    @(\"83\u00d721=\", \"18\u00d750=\"),\n    @(\"16\u00d743=\", \"40\u00d728=\"),\n    @(\"56\u00d748=\", \"25\u00d746=\"),\n    @(\"16\u00d748=\", \"82\u00d720=\"),\n    @(\"39\u00d730=\", \"41\u00d797=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $null, $wdReplaceAll) | Out-Null\n}\n"}
